# Updates cryptos list values (price/volume) and fixes a row-order swap
# between TheSandbox and InternetComputer(DFINITY), per the commit diff.
# Price cells in column D are numeric-looking text (e.g. '1.004', '0.00001101')
# that must remain stored as text, so they are written with a leading quote
# (quote-prefix) to stop Excel's automatic number/date conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.431.76"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "'1.798.87"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'338.27"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.3802"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").Value = "'48.59"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'1.203"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").Value = "'0.07502"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'22.13"
$ws.Range("E13").Value = "  +7.99%  "
$ws.Range("D14").Value = "'6.477"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "'1.795.62"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").Value = "'7.093"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "'0.06662"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'84.71"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "'6.525"
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("D22").Value = "'17.35"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("D23").Value = "'27.403.18"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").Value = "'12.55"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").Value = "'2.430"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").Value = "'2.565"
$ws.Range("E26").Value = "  +5.65%  "
$ws.Range("D27").Value = "'1.498"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'21.46"
$ws.Range("E28").Value = "  +9.36%  "
$ws.Range("D29").Value = "'151.92"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'1.998.70"
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("D31").Value = "'134.03"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").Value = "'4.056"
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("D33").Value = "'6.120"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").Value = "'0.08693"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'13.28"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").Value = "'1.640"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.457"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "'0.6915"
$ws.Range("E38").Value = "  +10.28%  "
$ws.Range("D39").Value = "'8.896"
$ws.Range("E39").Value = "  +4.36%  "
$ws.Range("D40").Value = "'0.06373"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").Value = "'0.2205"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").Value = "'0.02338"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "'1.277"
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("D44").Value = "'14.39"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "'0.6445"
$ws.Range("E45").Value = "  +5.75%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "'3.867"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "'2.136"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").Value = "'130.47"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "'0.07209"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'79.90"
$ws.Range("E51").Value = "  +2.67%  "
